# Add time frame selector with historical data support
#
# On the YouTube sheet: the "This Week" summary totals (row 4) double, and
# the "Top Videos" table (rows 9-17) shifts its data down by one row - each
# video's stats now additionally appear one row lower than before, as the
# newly added historical/comparison period data pushes the list down.
#
# NOTE: cell text such as "$114.60" or "6.8%" looks numeric to Excel's
# smart-entry parser and would otherwise be silently converted into a
# formatted number instead of staying literal text. A leading apostrophe
# (here embedded in a single-quoted PowerShell string, so it is never
# treated as $-variable interpolation) forces text entry, exactly like a
# user typing ' before the value in the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('YouTube')

# ---- Row 4: summary totals (doubled) ----
$ws.Range('A4').Value2 = 44208
$ws.Range('B4').Value2 = 780
$ws.Range('C4').Value = '''$114.60'
$ws.Range('F4').Value2 = 38

# ---- Top Videos table (rows 9-17): shifted down by one row ----

# Row 9
$ws.Range('A9').Value2 = 'I Built My Entire Design System in 4 Hours With AI'
$ws.Range('B9').Value2 = 9874
$ws.Range('C9').Value2 = 227
$ws.Range('D9').Value = '''6.8%'
$ws.Range('E9').Value = '''$23.75'

# Row 10
$ws.Range('A10').Value2 = 'Stop using V0 and Lovable to prototype, use Claude'
$ws.Range('B10').Value2 = 2885
$ws.Range('C10').Value2 = 40
$ws.Range('D10').Value = '''5.5%'
$ws.Range('E10').Value = '''$11.36'

# Row 11
$ws.Range('A11').Value2 = 'Stop using V0 and Lovable to prototype, use Claude'
$ws.Range('B11').Value2 = 2885
$ws.Range('C11').Value2 = 40
$ws.Range('D11').Value = '''5.5%'
$ws.Range('E11').Value = '''$11.36'

# Row 12
$ws.Range('A12').Value2 = 'How I Build a Component Library in 2 days (Figma t'
$ws.Range('B12').Value2 = 2352
$ws.Range('C12').Value2 = 58
$ws.Range('D12').Value = '''5.3%'
$ws.Range('E12').Value = '''$12.74'

# Row 13
$ws.Range('A13').Value2 = 'How I Build a Component Library in 2 days (Figma t'
$ws.Range('B13').Value2 = 2352
$ws.Range('C13').Value2 = 58
$ws.Range('D13').Value = '''5.3%'
$ws.Range('E13').Value = '''$12.74'
$ws.Range('F13').Value2 = 'No'

# Row 14
$ws.Range('A14').Value2 = 'Stop Wasting Dev Time on Frontend: Figma to Code i'
$ws.Range('B14').Value2 = 1247
$ws.Range('C14').Value2 = 16
$ws.Range('D14').Value = '''11.4%'
$ws.Range('E14').Value = '''$1.76'

# Row 15
$ws.Range('A15').Value2 = 'Stop Wasting Dev Time on Frontend: Figma to Code i'
$ws.Range('B15').Value2 = 1247
$ws.Range('C15').Value2 = 16
$ws.Range('D15').Value = '''11.4%'
$ws.Range('E15').Value = '''$1.76'

# Row 16
$ws.Range('A16').Value2 = 'Claude Code + Cursor + GitHub: The New AI environm'
$ws.Range('B16').Value2 = 999
$ws.Range('C16').Value2 = 20
$ws.Range('D16').Value = '''8.1%'
$ws.Range('E16').Value = '''$4.72'

# Row 17
$ws.Range('A17').Value2 = 'Claude Code + Cursor + GitHub: The New AI environm'
$ws.Range('B17').Value2 = 999
$ws.Range('C17').Value2 = 20
$ws.Range('D17').Value = '''8.1%'
$ws.Range('E17').Value = '''$4.72'
